$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 0)
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3 (Control 6)
$ws.Range("D3").Value = 0.00001164838593814021
$ws.Range("E3").Value = 0.00001164838593814021

# Row 4 (Control 9)
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = 0.00000000003840877913516637
$ws.Range("E4").Value = 0.00000000003840877913516637

# Row 5 (Control 24)
$ws.Range("D5").Value = 0.9924652230832347
$ws.Range("E5").Value = 0.9924652230832347

# Row 6 (Control 32)
$ws.Range("D6").Value = 0.9994793927500215
$ws.Range("E6").Value = 0.9994793927500215

# Row 7 (MDD 41)
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = 0.9889500467406946
$ws.Range("E7").Value = 0.01104995325930536

# Row 9 (MDD 15)
$ws.Range("D9").Value = 0.9999999999996232
$ws.Range("E9").Value = 0.0000000000003768096945577781

# Row 11 (MDD 33)
$ws.Range("D11").Value = 0.9999054824794873
$ws.Range("E11").Value = 0.00009451752051270468
$ws.Range("F11").Value = 5.133747577667236
$ws.Range("G11").Value = 0.7
